# Update the example-data-properties workbook:
#  - refresh the UBID values in column E (rows 2-15) with newly generated UBIDs
#  - bump a couple of row heights to fit the new values
#  - move the selection/scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New UBID values (column E, rows 2-15)
$ws.Range("E2").Value  = "86HJQCC9+5JJ-2-3-2-3"
$ws.Range("E3").Value  = "86HJX5QV+FJ3-2-3-2-2"
$ws.Range("E4").Value  = "86HJQ8Q5+R6V-1-2-1-1"
$ws.Range("E5").Value  = "86HJX6JP+H99-1-1-1-2"
$ws.Range("E6").Value  = "86HJQ76M+883-1-2-1-1"
$ws.Range("E7").Value  = "86HJW825+V3M-2-2-3-1"
$ws.Range("E8").Value  = "86HJX6GX+F4G-2-4-2-3"
$ws.Range("E9").Value  = "86HJX66G+P7C-2-3-2-3"
$ws.Range("E10").Value = "86HJM8JW+XMV-1-4-1-3"
$ws.Range("E11").Value = "86HJPCWJ+R59-1-5-2-4"
$ws.Range("E12").Value = "86HJR7QR+98F-2-1-1-1"
$ws.Range("E13").Value = "86HJQ9R3+FHW-1-2-0-3"
$ws.Range("E14").Value = "86HJW5RW+VGV-1-2-2-2"
$ws.Range("E15").Value = "86HJX838+8M7-1-3-1-2"

# Row height tweaks to accommodate the new text
$ws.Rows.Item(1).RowHeight = 26
$ws.Rows.Item(14).RowHeight = 32
$ws.Rows.Item(15).RowHeight = 32

# Reset the scroll position (was parked at O1) and move the selection
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E8").Select()
